$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.Execute("Responsive", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ins = $rng.Duplicate
$ins.Collapse(1)
$ins.InsertBefore("- ")
$ins.Font.Bold = $true
